# Scheduled runner update: refreshes cached Universalis market-price
# columns (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the leves whose
# listings moved since the last sync. Values are plain cached numbers (no
# formulas in this workbook), so each changed cell is written directly.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 833.3333
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 833.3333
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 833.3333
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1485.3333
$ws.Range("H33").Value = 1381.037
$ws.Range("I33").Value = 348.2353
$ws.Range("J33").Value = 3136.8
$ws.Range("K33").Value = 348.2353
$ws.Range("L33").Value = 3136.8
$ws.Range("M33").Value = -119.2353
$ws.Range("N33").Value = -3594.8
$ws.Range("H37").Value = 116
$ws.Range("J37").Value = 116
$ws.Range("L37").Value = 348
$ws.Range("N37").Value = -600
$ws.Range("H43").Value = 1934.2858
$ws.Range("I43").Value = 2497.625
$ws.Range("J43").Value = 1183.1666
$ws.Range("K43").Value = 2497.625
$ws.Range("L43").Value = 1183.1666
$ws.Range("M43").Value = -2428.625
$ws.Range("N43").Value = -1321.1666
$ws.Range("H64").Value = 37224.863
$ws.Range("I64").Value = 79489.234
$ws.Range("J64").Value = 2885.0625
$ws.Range("K64").Value = 79489.234
$ws.Range("L64").Value = 2885.0625
$ws.Range("M64").Value = -79241.234
$ws.Range("N64").Value = -3381.0625
$ws.Range("H67").Value = 37224.863
$ws.Range("I67").Value = 79489.234
$ws.Range("J67").Value = 2885.0625
$ws.Range("K67").Value = 79489.234
$ws.Range("L67").Value = 2885.0625
$ws.Range("M67").Value = -78631.234
$ws.Range("N67").Value = -4601.0625
$ws.Range("H69").Value = 3682.8572
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 3156
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 9468
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -11216
$ws.Range("H72").Value = 3682.8572
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 3156
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 28404
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -37140
$ws.Range("H74").Value = 3047.3157
$ws.Range("I74").Value = 2777.7778
$ws.Range("J74").Value = 3289.9
$ws.Range("K74").Value = 2777.7778
$ws.Range("L74").Value = 3289.9
$ws.Range("M74").Value = -1841.7778
$ws.Range("N74").Value = -5161.9
$ws.Range("H77").Value = 3047.3157
$ws.Range("I77").Value = 2777.7778
$ws.Range("J77").Value = 3289.9
$ws.Range("K77").Value = 13888.889
$ws.Range("L77").Value = 16449.5
$ws.Range("M77").Value = -9208.888999999999
$ws.Range("N77").Value = -25809.5
$ws.Range("H80").Value = 42737.125
$ws.Range("I80").Value = 998.5
$ws.Range("J80").Value = 63606.438
$ws.Range("K80").Value = 2995.5
$ws.Range("L80").Value = 190819.314
$ws.Range("M80").Value = -1997.5
$ws.Range("N80").Value = -192815.314
$ws.Range("H83").Value = 42737.125
$ws.Range("I83").Value = 998.5
$ws.Range("J83").Value = 63606.438
$ws.Range("K83").Value = 8986.5
$ws.Range("L83").Value = 572457.942
$ws.Range("M83").Value = -3994.5
$ws.Range("N83").Value = -582441.942

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 71324.87
$ws.Range("I45").Value = 92748.55
$ws.Range("J45").Value = 12409.75
$ws.Range("K45").Value = 92748.55
$ws.Range("L45").Value = 12409.75
$ws.Range("M45").Value = -92371.55
$ws.Range("N45").Value = -13163.75
$ws.Range("H122").Value = 1572.3549
$ws.Range("I122").Value = 1488.2727
$ws.Range("K122").Value = 4464.8181
$ws.Range("M122").Value = -2014.8181
$ws.Range("H132").Value = 10650.242
$ws.Range("I132").Value = 12681.2705
$ws.Range("K132").Value = 38043.8115
$ws.Range("M132").Value = -35513.8115

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 75211.734
$ws.Range("I86").Value = 101831.45
$ws.Range("J86").Value = 2007.5
$ws.Range("K86").Value = 101831.45
$ws.Range("L86").Value = 2007.5
$ws.Range("M86").Value = -100708.45
$ws.Range("N86").Value = -4253.5
$ws.Range("H89").Value = 75211.734
$ws.Range("I89").Value = 101831.45
$ws.Range("J89").Value = 2007.5
$ws.Range("K89").Value = 509157.25
$ws.Range("L89").Value = 10037.5
$ws.Range("M89").Value = -503541.25
$ws.Range("N89").Value = -21269.5
$ws.Range("H99").Value = 2127.8965
$ws.Range("I99").Value = 1648.5454
$ws.Range("K99").Value = 1648.5454
$ws.Range("M99").Value = -150.5454
$ws.Range("H105").Value = 62110
$ws.Range("I105").Value = 39976.46
$ws.Range("J105").Value = 144320.28
$ws.Range("K105").Value = 39976.46
$ws.Range("L105").Value = 144320.28
$ws.Range("M105").Value = -38229.46
$ws.Range("N105").Value = -147814.28
$ws.Range("H107").Value = 55556390
$ws.Range("I107").Value = 62500816
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 62500816
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -62498896
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 2107.5454
$ws.Range("I134").Value = 1499.8334
$ws.Range("K134").Value = 4499.5002
$ws.Range("M134").Value = -1964.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31821.896
$ws.Range("I31").Value = 627.5789
$ws.Range("J31").Value = 52259.55
$ws.Range("K31").Value = 627.5789
$ws.Range("L31").Value = 52259.55
$ws.Range("M31").Value = -332.5789
$ws.Range("N31").Value = -52849.55
$ws.Range("H33").Value = 3140.4285
$ws.Range("I33").Value = 3140.4285
$ws.Range("K33").Value = 3140.4285
$ws.Range("M33").Value = -2761.4285
$ws.Range("H34").Value = 31821.896
$ws.Range("I34").Value = 627.5789
$ws.Range("J34").Value = 52259.55
$ws.Range("K34").Value = 627.5789
$ws.Range("L34").Value = 52259.55
$ws.Range("M34").Value = -425.5789
$ws.Range("N34").Value = -52663.55
$ws.Range("H99").Value = 9062.733
$ws.Range("I99").Value = 2098
$ws.Range("J99").Value = 11595.363
$ws.Range("K99").Value = 2098
$ws.Range("L99").Value = 11595.363
$ws.Range("M99").Value = -600
$ws.Range("N99").Value = -14591.363
$ws.Range("H126").Value = 9062.733
$ws.Range("I126").Value = 2098
$ws.Range("J126").Value = 11595.363
$ws.Range("K126").Value = 6294
$ws.Range("L126").Value = 34786.089
$ws.Range("M126").Value = -3824
$ws.Range("N126").Value = -39726.089

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1692.2333
$ws.Range("I122").Value = 1498.1666
$ws.Range("K122").Value = 4494.4998
$ws.Range("M122").Value = -2044.4998
$ws.Range("H126").Value = 1927.7391
$ws.Range("I126").Value = 1866.7368
$ws.Range("J126").Value = 2217.5
$ws.Range("K126").Value = 5600.2104
$ws.Range("L126").Value = 6652.5
$ws.Range("M126").Value = -3130.2104
$ws.Range("N126").Value = -11592.5
$ws.Range("H132").Value = 2446.4
$ws.Range("I132").Value = 1891.6487
$ws.Range("J132").Value = 5012.125
$ws.Range("K132").Value = 5674.9461
$ws.Range("L132").Value = 15036.375
$ws.Range("M132").Value = -3144.9461
$ws.Range("N132").Value = -20096.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 393592.12
$ws.Range("I46").Value = 279
$ws.Range("J46").Value = 601816.7
$ws.Range("K46").Value = 279
$ws.Range("L46").Value = 601816.7
$ws.Range("M46").Value = -91
$ws.Range("N46").Value = -602192.7
$ws.Range("H55").Value = 253168.95
$ws.Range("I55").Value = 569037.1
$ws.Range("J55").Value = 474.4
$ws.Range("K55").Value = 569037.1
$ws.Range("L55").Value = 474.4
$ws.Range("M55").Value = -568864.1
$ws.Range("N55").Value = -820.4
$ws.Range("H82").Value = 2075.5
$ws.Range("I82").Value = 1652.4
$ws.Range("J82").Value = 2310.5557
$ws.Range("K82").Value = 1652.4
$ws.Range("L82").Value = 2310.5557
$ws.Range("M82").Value = -1291.4
$ws.Range("N82").Value = -3032.5557
$ws.Range("H85").Value = 2075.5
$ws.Range("I85").Value = 1652.4
$ws.Range("J85").Value = 2310.5557
$ws.Range("K85").Value = 1652.4
$ws.Range("L85").Value = 2310.5557
$ws.Range("M85").Value = -404.4000000000001
$ws.Range("N85").Value = -4806.5557
$ws.Range("H122").Value = 2499.7693
$ws.Range("I122").Value = 2401.6365
$ws.Range("J122").Value = 3039.5
$ws.Range("K122").Value = 7204.9095
$ws.Range("L122").Value = 9118.5
$ws.Range("M122").Value = -4754.9095
$ws.Range("N122").Value = -14018.5
$ws.Range("H132").Value = 5105.857
$ws.Range("I132").Value = 5201.6
$ws.Range("K132").Value = 15604.8
$ws.Range("M132").Value = -13074.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 400437.6
$ws.Range("I81").Value = 333829.66
$ws.Range("J81").Value = 500349.5
$ws.Range("K81").Value = 667659.3199999999
$ws.Range("L81").Value = 1000699
$ws.Range("M81").Value = -666598.3199999999
$ws.Range("N81").Value = -1002821
$ws.Range("H84").Value = 400437.6
$ws.Range("I84").Value = 333829.66
$ws.Range("J84").Value = 500349.5
$ws.Range("K84").Value = 3338296.6
$ws.Range("L84").Value = 5003495
$ws.Range("M84").Value = -3332992.6
$ws.Range("N84").Value = -5014103
$ws.Range("H122").Value = 1427.2188
$ws.Range("I122").Value = 1470.3928
$ws.Range("K122").Value = 4411.178400000001
$ws.Range("M122").Value = -1961.178400000001
